$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 560. This shifts existing rows 560:597 down to 561:598,
# preserving all their data/formatting, and grows the sheet from R597 to R598.
$ws.Rows("560:560").Insert()

# Populate the newly inserted row 560 with the new week's data.
$ws.Range("A560").Value = 6
$ws.Range("B560").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C560").Value = "Metropolitana"
$ws.Range("D560").Value = 44826
$ws.Range("E560").Value = 13
$ws.Range("F560").Value = 100112044
$ws.Range("G560").Value = "Perejil"
$ws.Range("H560").Value = "Sin especificar"
$ws.Range("I560").Value = "Primera"
$ws.Range("J560").Value = 280
$ws.Range("K560").Value = 11000
$ws.Range("L560").Value = 12000
$ws.Range("M560").Value = 11393
$ws.Range("N560").Value = '$/docena de atados'
$ws.Range("O560").Value = "Región Metropolitana"
$ws.Range("P560").Value = 3798
$ws.Range("Q560").Value = 3
$ws.Range("R560").Value = "Hortaliza"
